$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): F2 417->418, F3 2610->2625, F4 118->119
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 418
$ws1.Range("F3").Value = 2625
$ws1.Range("F4").Value = 119

# Sheet "全部类型" (sheet4): F2 417->418, F7 2610->2625, F8 118->119
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 418
$ws4.Range("F7").Value = 2625
$ws4.Range("F8").Value = 119
